$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F values
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F2").Value = 4705
$wsExh.Range("F4").Value = 72
$wsExh.Range("F6").Value = 242
$wsExh.Range("F9").Value = 190
$wsExh.Range("F10").Value = 1820
$wsExh.Range("F11").Value = 328
$wsExh.Range("F12").Value = 4207
$wsExh.Range("F13").Value = 52
$wsExh.Range("F14").Value = 289

# Sheet "全部类型" (sheet4): update column F values (aggregated rows)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4705
$wsAll.Range("F5").Value = 72
$wsAll.Range("F8").Value = 242
$wsAll.Range("F11").Value = 190
$wsAll.Range("F14").Value = 1820
$wsAll.Range("F15").Value = 328
$wsAll.Range("F16").Value = 4207
$wsAll.Range("F17").Value = 52
$wsAll.Range("F18").Value = 289
